$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-25 09:15:45"
$wsZhCn.Range("H3").Value = "2016-03-25 09:16:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-25 09:15:49"
$wsDeDe.Range("H3").Value = "2016-03-25 09:16:30"
